$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "run"
$ws.Range("B1").Value = "pheno"
$ws.Range("C1").Value = "biorep"

# Data rows: run id, pheno (WT/KO), biorep (b_1..b_6)
$data = @(
    @("1WT_20_2h_n3_1", "WT", "b_1"),
    @("1WT_20_2h_n3_2", "WT", "b_1"),
    @("1WT_20_2h_n3_3", "WT", "b_1"),
    @("1WT_20_2h_n4_1", "WT", "b_2"),
    @("1WT_20_2h_n4_2", "WT", "b_2"),
    @("1WT_20_2h_n4_3", "WT", "b_2"),
    @("1WT_20_2h_n5_1", "WT", "b_3"),
    @("1WT_20_2h_n5_2", "WT", "b_3"),
    @("1WT_20_2h_n5_3", "WT", "b_3"),
    @("3D8_20_2h_n3_1", "KO", "b_4"),
    @("3D8_20_2h_n3_2", "KO", "b_4"),
    @("3D8_20_2h_n3_3", "KO", "b_4"),
    @("3D8_20_2h_n4_1", "KO", "b_5"),
    @("3D8_20_2h_n4_2", "KO", "b_5"),
    @("3D8_20_2h_n4_3", "KO", "b_5"),
    @("3D8_20_2h_n5_1", "KO", "b_6"),
    @("3D8_20_2h_n5_2", "KO", "b_6"),
    @("3D8_20_2h_n5_3", "KO", "b_6")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Column width: move bestFit custom width from column C to column A.
# (15.33 is the closest input that this engine's pixel-quantization maps to
# a stored width nearest the source value of 16.1640625.)
$ws.Columns.Item(1).ColumnWidth = 15.33
